$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 102
$ws.Range("E2").Value = 86
$ws.Range("F2").Value = 0.8431372549019608
$ws.Range("G2").Value = 0.8431372549019608
$ws.Range("H2").Value = 0.09877977698439171
$ws.Range("I2").Value = 0.0832849100064479
$ws.Range("J2").Value = 461132.027355649
$ws.Range("K2").Value = 167952.0138788245
$ws.Range("M2").Value = 167952.0138788245
$ws.Range("N2").Value = 629084.0412344737
$ws.Range("O2").Value = 10240780.2488
$ws.Range("P2").Value = 9833039.318699999
$ws.Range("Q2").Value = 0.01640031421419328
$ws.Range("R2").Value = 0.01708037651791155

# Row 3
$ws.Range("C3").Value = 106
$ws.Range("D3").Value = 106
$ws.Range("E3").Value = 91
$ws.Range("F3").Value = 0.8584905660377359
$ws.Range("G3").Value = 0.8584905660377359
$ws.Range("H3").Value = 0.09496116222475419
$ws.Range("I3").Value = 0.08152326190993048
$ws.Range("J3").Value = 486228.6741479071
$ws.Range("K3").Value = 177901.9172809835
$ws.Range("M3").Value = 177901.9172809835
$ws.Range("N3").Value = 664130.5914288907
$ws.Range("O3").Value = 10797211.559364
$ws.Range("P3").Value = 10389838.401361
$ws.Range("Q3").Value = 0.01647665383815658
$ws.Range("R3").Value = 0.01712268376163383

# Row 4
$ws.Range("C4").Value = 109
$ws.Range("D4").Value = 109
$ws.Range("E4").Value = 93
$ws.Range("F4").Value = 0.8532110091743119
$ws.Range("G4").Value = 0.8532110091743119
$ws.Range("H4").Value = 0.09300393342188501
$ws.Range("I4").Value = 0.07935197989206702
$ws.Range("J4").Value = 512838.5062540149
$ws.Range("K4").Value = 184120.4607402484
$ws.Range("M4").Value = 184120.4607402484
$ws.Range("N4").Value = 696958.9669942633
$ws.Range("O4").Value = 11361855.35814492
$ws.Range("P4").Value = 10953511.00540183
$ws.Range("Q4").Value = 0.01620514035221006
$ws.Range("R4").Value = 0.01680926423038673

# Row 5
$ws.Range("C5").Value = 110
$ws.Range("D5").Value = 109
$ws.Range("E5").Value = 93
$ws.Range("F5").Value = 0.8532110091743119
$ws.Range("G5").Value = 0.8454545454545455
$ws.Range("H5").Value = 0.09289780249835385
$ws.Range("I5").Value = 0.0785408693849719
$ws.Range("J5").Value = 523961.5408676272
$ws.Range("K5").Value = 187108.0142754518
$ws.Range("M5").Value = 187108.0142754518
$ws.Range("N5").Value = 711069.5551430788
$ws.Range("O5").Value = 11517660.85818927
$ws.Range("P5").Value = 11106966.17486389
$ws.Range("Q5").Value = 0.01624531374722798
$ws.Range("R5").Value = 0.01684600559051803

# Row 6
$ws.Range("C6").Value = 112
$ws.Range("D6").Value = 112
$ws.Range("E6").Value = 95
$ws.Range("F6").Value = 0.8482142857142857
$ws.Range("G6").Value = 0.8482142857142857
$ws.Range("H6").Value = 0.09190722738373355
$ws.Range("I6").Value = 0.07795702322727401
$ws.Range("J6").Value = 548708.848341326
$ws.Range("K6").Value = 196420.4853275503
$ws.Range("M6").Value = 196420.4853275503
$ws.Range("N6").Value = 745129.3336688762
$ws.Range("O6").Value = 12112094.92793495
$ws.Range("P6").Value = 11697629.4041098
$ws.Range("Q6").Value = 0.01621688787086141
$ws.Range("R6").Value = 0.01679147787487101

# Row 7
$ws.Range("D7").Value = 102
$ws.Range("F7").Value = 0.8529411764705882
$ws.Range("H7").Value = 0.09823822539669926
$ws.Range("I7").Value = 0.08379142754424349
$ws.Range("J7").Value = 461915.3705097084
$ws.Range("K7").Value = 168343.6854558542
$ws.Range("M7").Value = 168343.6854558542
$ws.Range("N7").Value = 630259.0559655628
$ws.Range("O7").Value = 10165519.7688
$ws.Range("P7").Value = 9757778.838699998
$ws.Range("Q7").Value = 0.01656026344786958
$ws.Range("R7").Value = 0.01725225466150062

# Row 8
$ws.Range("C8").Value = 106
$ws.Range("D8").Value = 106
$ws.Range("E8").Value = 91
$ws.Range("F8").Value = 0.8584905660377359
$ws.Range("G8").Value = 0.8584905660377359
$ws.Range("H8").Value = 0.09607362778140782
$ws.Range("I8").Value = 0.08247830309535956
$ws.Range("J8").Value = 493240.0544632261
$ws.Range("K8").Value = 181407.607438643
$ws.Range("M8").Value = 181407.607438643
$ws.Range("N8").Value = 674647.6619018689
$ws.Range("O8").Value = 10801469.936564
$ws.Range("P8").Value = 10394096.778561
$ws.Range("Q8").Value = 0.01679471484011274
$ws.Range("R8").Value = 0.01745294577329862

# Row 9
$ws.Range("C9").Value = 109
$ws.Range("D9").Value = 109
$ws.Range("E9").Value = 94
$ws.Range("F9").Value = 0.8623853211009175
$ws.Range("G9").Value = 0.8623853211009175
$ws.Range("H9").Value = 0.09422337361318721
$ws.Range("I9").Value = 0.08125685430862017
$ws.Range("J9").Value = 526514.418027284
$ws.Range("K9").Value = 190958.4166268829
$ws.Range("M9").Value = 190958.4166268829
$ws.Range("N9").Value = 717472.8346541669
$ws.Range("O9").Value = 11345286.94306092
$ws.Range("P9").Value = 10936942.59031783
$ws.Range("Q9").Value = 0.01683151934237134
$ws.Range("R9").Value = 0.017459945048622

# Row 10
$ws.Range("C10").Value = 110
$ws.Range("D10").Value = 110
$ws.Range("E10").Value = 95
$ws.Range("F10").Value = 0.8636363636363636
$ws.Range("G10").Value = 0.8636363636363636
$ws.Range("H10").Value = 0.09389805491180749
$ws.Range("I10").Value = 0.08109377469656102
$ws.Range("J10").Value = 548242.7938357895
$ws.Range("K10").Value = 199248.6407595329
$ws.Range("M10").Value = 199248.6407595329
$ws.Range("N10").Value = 747491.4345953225
$ws.Range("O10").Value = 11786101.98315275
$ws.Range("P10").Value = 11375407.29982737
$ws.Range("Q10").Value = 0.016905389164657
$ws.Range("R10").Value = 0.01751573684421451

# Row 11
$ws.Range("C11").Value = 112
$ws.Range("D11").Value = 112
$ws.Range("E11").Value = 97
$ws.Range("F11").Value = 0.8660714285714286
$ws.Range("G11").Value = 0.8660714285714286
$ws.Range("H11").Value = 0.09400854608248253
$ws.Range("I11").Value = 0.08141811580357862
$ws.Range("J11").Value = 574500.8661924924
$ws.Range("K11").Value = 209316.4942531334
$ws.Range("M11").Value = 209316.4942531334
$ws.Range("N11").Value = 783817.3604456257
$ws.Range("O11").Value = 12106775.67824733
$ws.Range("P11").Value = 11692310.15442219
$ws.Range("Q11").Value = 0.01728920232900819
$ws.Range("R11").Value = 0.01790206481770133

# Row 12
$ws.Range("D12").Value = 102
$ws.Range("E12").Value = 102
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 0.09127332586241999
$ws.Range("I12").Value = 0.09127332586241999
$ws.Range("J12").Value = 544219.5395527922
$ws.Range("K12").Value = 209495.769977396
$ws.Range("M12").Value = 209495.769977396
$ws.Range("N12").Value = 753715.3095301883
$ws.Range("O12").Value = 10316742.6388
$ws.Range("P12").Value = 9909001.708699998
$ws.Range("Q12").Value = 0.02030638713323217
$ws.Range("R12").Value = 0.02114196526916137

# Row 13
$ws.Range("C13").Value = 106
$ws.Range("D13").Value = 106
$ws.Range("E13").Value = 106
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 0.09717788792471758
$ws.Range("I13").Value = 0.09717788792471758
$ws.Range("J13").Value = 654630.7013274725
$ws.Range("K13").Value = 262102.9308707663
$ws.Range("M13").Value = 262102.9308707663
$ws.Range("N13").Value = 916733.6321982386
$ws.Range("O13").Value = 10760730.661864
$ws.Range("P13").Value = 10353357.503861
$ws.Range("Q13").Value = 0.02435735444988493
$ws.Range("R13").Value = 0.02531574233508523

# Row 14
$ws.Range("C14").Value = 108
$ws.Range("D14").Value = 108
$ws.Range("E14").Value = 108
$ws.Range("H14").Value = 0.09963069641345117
$ws.Range("I14").Value = 0.09963069641345117
$ws.Range("J14").Value = 734817.7333005213
$ws.Range("K14").Value = 295110.0742635016
$ws.Range("M14").Value = 295110.0742635016
$ws.Range("N14").Value = 1029927.807564023
$ws.Range("O14").Value = 11223740.87971992
$ws.Range("P14").Value = 10815396.52697683
$ws.Range("Q14").Value = 0.02629337913500242
$ws.Range("R14").Value = 0.02728610768245148

# Row 15
$ws.Range("C15").Value = 109
$ws.Range("D15").Value = 109
$ws.Range("E15").Value = 109
$ws.Range("H15").Value = 0.1015414746775767
$ws.Range("I15").Value = 0.1015414746775767
$ws.Range("J15").Value = 781388.0653614923
$ws.Range("K15").Value = 315821.2765223843
$ws.Range("M15").Value = 315821.2765223843
$ws.Range("N15").Value = 1097209.341883876
$ws.Range("O15").Value = 11501137.06941152
$ws.Range("P15").Value = 11090442.38608613
$ws.Range("Q15").Value = 0.02746000457314296
$ws.Range("R15").Value = 0.02847688717256292

# Row 16
$ws.Range("C16").Value = 110
$ws.Range("D16").Value = 110
$ws.Range("E16").Value = 110
$ws.Range("H16").Value = 0.1022024053100728
$ws.Range("I16").Value = 0.1022024053100728
$ws.Range("J16").Value = 827577.8738163244
$ws.Range("K16").Value = 335854.9980650494
$ws.Range("M16").Value = 335854.9980650494
$ws.Range("N16").Value = 1163432.871881374
$ws.Range("O16").Value = 11926849.94239386
$ws.Range("P16").Value = 11512384.41856872
$ws.Range("Q16").Value = 0.02815957270253366
$ws.Range("R16").Value = 0.02917336546921916
